# Rewrites the LM Week document body to the new "Week 12" / team-synthesis
# version described by the commit, and shrinks the page margins to match.
$d = $word.ActiveDocument

# Build the replacement OOXML for the document body as a literal (non-interpolating)
# here-string, then push it in through Range.InsertXML -- this restates the whole
# body (paragraphs, runs, run formatting, and both tables) in one shot, which is far
# less error-prone than hundreds of piecemeal Find/Replace + paragraph-insert calls
# for a change this broad. Range.Content stops short of the trailing section mark,
# so the existing sectPr (and its rsids) survive untouched.
$newBodyXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:b/>
          <w:sz w:val="32"/>
        </w:rPr>
        <w:t>NotebookLM Discussion Post - Week 12</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:b/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t>TCE 486/586A</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t>Spring 2026 (Edwards)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t>10 points possible</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:b/>
          <w:color w:val="0F4761"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t>1. What is this all about?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t>These discussion posts have two purposes: (1) To offer you a chance to reflect on the readings and discussions you had with classmates about them; and (2) To provide you with practice using genAI tools, specifically NotebookLM, to engage with readings. This anticipates work you will likely do with your own students as well as providing you with tools to improve your own reading comprehension.</w:t>
        <w:br/>
        <w:br/>
        <w:t>This Week's Prompts:</w:t>
        <w:br/>
        <w:br/>
        <w:t>Boaler, Chapter 7: From Tracking to Growth Mindset Grouping</w:t>
        <w:br/>
        <w:t>Boaler says ability grouping (tracking) delivers "fixed mindset messages"—but what if your school REQUIRES leveled classes? Can you create growth mindset within a tracked system, or does the structure itself undermine everything you're trying to do?</w:t>
        <w:br/>
        <w:br/>
        <w:t>Ambitious Science Teaching, Chapter 4: Talk in Science</w:t>
        <w:br/>
        <w:t>If productive science talk means students "build on each other's ideas"—what happens when one student dominates the discussion or when the quiet kids never speak? Do we intervene and risk shutting down organic conversation, or stay silent and accept inequity?</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:b/>
          <w:color w:val="0F4761"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t>2. What do I need to do?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t>Answer the philosophical prompt for the chapter YOU deep-read this week, then synthesize what you learned from your team.</w:t>
        <w:br/>
        <w:br/>
        <w:t>Part 1: Chapter Deep-Dive (4 points)</w:t>
        <w:br/>
        <w:br/>
        <w:t>Choose the prompt for YOUR chapter:</w:t>
        <w:br/>
        <w:br/>
        <w:t>Prompt A (Boaler Ch 7 - From Tracking to Growth Mindset Grouping):</w:t>
        <w:br/>
        <w:t>Boaler says ability grouping (tracking) delivers "fixed mindset messages"—but what if your school REQUIRES leveled classes? Can you create growth mindset within a tracked system, or does the structure itself undermine everything you're trying to do?</w:t>
        <w:br/>
        <w:br/>
        <w:t>Prompt B (Ambitious Science Ch 4 - Talk in Science):</w:t>
        <w:br/>
        <w:t>If productive science talk means students "build on each other's ideas"—what happens when one student dominates the discussion or when the quiet kids never speak? Do we intervene and risk shutting down organic conversation, or stay silent and accept inequity?</w:t>
        <w:br/>
        <w:br/>
        <w:t>After responding to your prompt, add:</w:t>
        <w:br/>
        <w:t>• What NotebookLM tool(s) did you use this week? (Audio Overview, Study Guide, Briefing Doc, FAQ, etc.)</w:t>
        <w:br/>
        <w:t>• What's one teaching idea from YOUR chapter that you want to try? Be specific!</w:t>
        <w:br/>
        <w:br/>
        <w:t>Part 2: Synthesis - Learning from Your Team (6 points)</w:t>
        <w:br/>
        <w:br/>
        <w:t>Human Synthesis (3 points):</w:t>
        <w:br/>
        <w:t>What did you learn from your teammates about the OTHER chapter (the one you didn't deep-read)? What's one key idea or teaching strategy from that chapter that stuck with you from Thursday's discussion?</w:t>
        <w:br/>
        <w:br/>
        <w:t>AI Synthesis (3 points):</w:t>
        <w:br/>
        <w:t>Go to NotebookLM and ask the AI a synthesis question that connects BOTH chapters. For example: "How does productive talk (AS Ch 4) work in heterogeneous vs tracked classrooms (Boaler Ch 7)?" or "What talk moves support equity in mixed-ability groups?"</w:t>
        <w:br/>
        <w:t>Include a screenshot of your question and the AI's response.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:b/>
          <w:color w:val="0F4761"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t>3. How will my work be assessed? (10 points)</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:tbl>
      <w:tblPr>
        <w:tblStyle w:val="LightGrid-Accent1"/>
        <w:tblW w:type="auto" w:w="0"/>
        <w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="3600"/>
        <w:gridCol w:w="3600"/>
        <w:gridCol w:w="3600"/>
      </w:tblGrid>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="3600"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:b/>
                <w:sz w:val="22"/>
              </w:rPr>
              <w:t>Criterion</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="3600"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:b/>
                <w:sz w:val="22"/>
              </w:rPr>
              <w:t>Points</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="3600"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:b/>
                <w:sz w:val="22"/>
              </w:rPr>
              <w:t>What I'm Looking For</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="3600"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="22"/>
              </w:rPr>
              <w:t>Chapter deep-dive</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="3600"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="22"/>
              </w:rPr>
              <w:t>4</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="3600"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="22"/>
              </w:rPr>
              <w:t>Thoughtful response to philosophical prompt + specific NotebookLM tools used + concrete teaching idea from YOUR chapter</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="3600"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="22"/>
              </w:rPr>
              <w:t>Human synthesis</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="3600"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="22"/>
              </w:rPr>
              <w:t>3</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="3600"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="22"/>
              </w:rPr>
              <w:t>Clear explanation of what you learned from teammates about the OTHER chapter with specific teaching strategy</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="3600"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="22"/>
              </w:rPr>
              <w:t>AI synthesis</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="3600"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="22"/>
              </w:rPr>
              <w:t>3</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="3600"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="22"/>
              </w:rPr>
              <w:t>Synthesis question connecting both chapters + screenshot of NotebookLM response</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
    <w:p/>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:b/>
          <w:color w:val="0F4761"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t>4. Submission</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t>Submit your work to Canvas.</w:t>
      </w:r>
    </w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$d.Content.InsertXML($newBodyXml) | Out-Null

# Shrink the page margins from 1440/1800 twips to 720 twips on every side (1440/1800
# twips == 72/90 points; 720 twips == 36 points -- Word COM reports margins in points).
$d.PageSetup.TopMargin = 36
$d.PageSetup.BottomMargin = 36
$d.PageSetup.LeftMargin = 36
$d.PageSetup.RightMargin = 36

Write-Output "Paragraphs: $($d.Paragraphs.Count); Tables: $($d.Tables.Count)"
